$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "# Past Due" (B) and "# Extensions" (C) counts for rows 2-61
$pastDue = @(
    0, 7, 0, 1, 10, 3, 3, 5, 1, 3, 2, 3, 0, 3, 0, 0, 0, 9, 1, 0, 0, 3, 3, 2, 0, 3, 7, 6, 2, 3, 1, 8, 3, 2, 0, 22, 2, 1, 1, 0, 2, 3, 2, 7, 8, 3, 3, 2, 5, 2, 3, 0, 1, 1, 0, 2, 0, 0, 2, 3
)

$extensions = @(
    2, 3, 1, 1, 3, 1, 0, 2, 0, 1, 0, 0, 1, 0, 1, 1, 1, 3, 1, 0, 2, 1, 1, 1, 2, 1, 4, 3, 2, 1, 1, 4, 1, 2, 2, 9, 0, 2, 0, 1, 1, 0, 2, 4, 3, 1, 1, 1, 4, 2, 1, 1, 0, 2, 0, 2, 2, 0, 1, 2
)

for ($i = 0; $i -lt $pastDue.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $pastDue[$i]
    $ws.Cells.Item($row, 3).Value = $extensions[$i]
}

